$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035885148950424
$ws.Range("D2").Value = 1.038311417493139
$ws.Range("E2").Value = 1.044491955831546
$ws.Range("F2").Value = 1.054574073263947
$ws.Range("I2").Value = 1.038838669487082
$ws.Range("J2").Value = 1.040996482544607
$ws.Range("K2").Value = 1.041099811562334
$ws.Range("L2").Value = 1.04726285444757
$ws.Range("M2").Value = 1.057316909141941
$ws.Range("N2").Value = 1.042474816953786

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036763236396237
$ws.Range("D3").Value = 1.03895692447348
$ws.Range("E3").Value = 1.045345701130501
$ws.Range("F3").Value = 1.055677028276694
$ws.Range("I3").Value = 1.03905861747575
$ws.Range("J3").Value = 1.041518664465915
$ws.Range("K3").Value = 1.041555740632425
$ws.Range("L3").Value = 1.047927737725493
$ws.Range("M3").Value = 1.058232394153206
$ws.Range("N3").Value = 1.042997740433317

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037331717464131
$ws.Range("D4").Value = 1.039374755357087
$ws.Range("E4").Value = 1.045898845965253
$ws.Range("F4").Value = 1.056391953457573
$ws.Range("I4").Value = 1.039199683201823
$ws.Range("J4").Value = 1.041856204620144
$ws.Range("K4").Value = 1.041850215496057
$ws.Range("L4").Value = 1.048358032575444
$ws.Range("M4").Value = 1.058825419729738
$ws.Range("N4").Value = 1.043335759933282

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037570777521864
$ws.Range("D5").Value = 1.039550444676601
$ws.Range("E5").Value = 1.046131557837729
$ws.Range("F5").Value = 1.056692803672369
$ws.Range("I5").Value = 1.039258686073731
$ws.Range("J5").Value = 1.041998022594004
$ws.Range("K5").Value = 1.041973882050207
$ws.Range("L5").Value = 1.048538944489139
$ws.Range("M5").Value = 1.059074881185133
$ws.Range("N5").Value = 1.043477779304931

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037610920865109
$ws.Range("D6").Value = 1.039579945592711
$ws.Range("E6").Value = 1.046170641104191
$ws.Range("F6").Value = 1.056743335054156
$ws.Range("I6").Value = 1.039268575237103
$ws.Range("J6").Value = 1.042021829507535
$ws.Range("K6").Value = 1.041994638517061
$ws.Range("L6").Value = 1.048569321301795
$ws.Range("M6").Value = 1.059116775832381
$ws.Range("N6").Value = 1.04350162002701

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037334911520282
$ws.Range("D7").Value = 1.039377102797652
$ws.Range("E7").Value = 1.045901954807148
$ws.Range("F7").Value = 1.056395972272433
$ws.Range("I7").Value = 1.039200472785223
$ws.Range("J7").Value = 1.041858099929755
$ws.Range("K7").Value = 1.041851868449322
$ws.Range("L7").Value = 1.048360449867854
$ws.Range("M7").Value = 1.058828752444055
$ws.Range("N7").Value = 1.04333765793445

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036181839951708
$ws.Range("D8").Value = 1.038529538780837
$ws.Range("E8").Value = 1.044780334334161
$ws.Range("F8").Value = 1.05494656477083
$ws.Range("I8").Value = 1.038913261484623
$ws.Range("J8").Value = 1.0411730275293
$ws.Range("K8").Value = 1.041254006502272
$ws.Range("L8").Value = 1.047487539570126
$ws.Range("M8").Value = 1.057626167291697
$ws.Range("N8").Value = 1.042651612652607

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034152327184176
$ws.Range("D9").Value = 1.037037193835469
$ws.Range("E9").Value = 1.042809424674364
$ws.Range("F9").Value = 1.052402064578012
$ws.Range("I9").Value = 1.038397571985315
$ws.Range("J9").Value = 1.03996323135781
$ws.Range("K9").Value = 1.040196391861646
$ws.Range("L9").Value = 1.045949948074749
$ws.Range("M9").Value = 1.055512048950903
$ws.Range("N9").Value = 1.041440098431788

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032800963763663
$ws.Range("D10").Value = 1.036043168919097
$ws.Range("E10").Value = 1.041499277488244
$ws.Range("F10").Value = 1.050712208884296
$ws.Range("I10").Value = 1.038047370699735
$ws.Range("J10").Value = 1.039155001635327
$ws.Range("K10").Value = 1.03948861326312
$ws.Range("L10").Value = 1.044925342173922
$ws.Range("M10").Value = 1.054106053071604
$ws.Range("N10").Value = 1.040630720930394

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032216212519283
$ws.Range("D11").Value = 1.035612969544607
$ws.Range("E11").Value = 1.040932884457755
$ws.Range("F11").Value = 1.049982031997249
$ws.Range("I11").Value = 1.037894217400124
$ws.Range("J11").Value = 1.038804638002702
$ws.Range("K11").Value = 1.039181508479518
$ws.Range("L11").Value = 1.044481796806557
$ws.Range("M11").Value = 1.053498064175034
$ws.Range("N11").Value = 1.040279859741223

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.031999070741217
$ws.Range("D12").Value = 1.035453208794644
$ws.Range("E12").Value = 1.040722638717341
$ws.Range("F12").Value = 1.049711044413691
$ws.Range("I12").Value = 1.037837102595132
$ws.Range("J12").Value = 1.038674438944481
$ws.Range("K12").Value = 1.039067342085986
$ws.Range("L12").Value = 1.04431706266077
$ws.Range("M12").Value = 1.053272353678993
$ws.Range("N12").Value = 1.040149475785406

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032045645618006
$ws.Range("D13").Value = 1.03548747643628
$ws.Range("E13").Value = 1.040767730884593
$ws.Range("F13").Value = 1.049769161627631
$ws.Range("I13").Value = 1.037849364170723
$ws.Range("J13").Value = 1.038702369712724
$ws.Range("K13").Value = 1.039091835403617
$ws.Range("L13").Value = 1.044352397848284
$ws.Range("M13").Value = 1.053320763735568
$ws.Range("N13").Value = 1.040177446218543

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032198262264237
$ws.Range("D14").Value = 1.035599762960251
$ws.Range("E14").Value = 1.04091550266005
$ws.Range("F14").Value = 1.049959627327874
$ws.Range("I14").Value = 1.037889500896812
$ws.Range("J14").Value = 1.038793876894095
$ws.Range("K14").Value = 1.039172073361984
$ws.Range("L14").Value = 1.044468179450719
$ws.Range("M14").Value = 1.053479404346976
$ws.Range("N14").Value = 1.040269083350607

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032292302514914
$ws.Range("D15").Value = 1.035668951000809
$ws.Range("E15").Value = 1.04100656803126
$ws.Range("F15").Value = 1.050077010386478
$ws.Range("I15").Value = 1.037914200416105
$ws.Range("J15").Value = 1.038850249776748
$ws.Range("K15").Value = 1.039221498188261
$ws.Range("L15").Value = 1.044539518777963
$ws.Range("M15").Value = 1.053577164477384
$ws.Range("N15").Value = 1.04032553628922

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03283978025531
$ws.Range("D16").Value = 1.036071724589235
$ws.Range("E16").Value = 1.041536886392221
$ws.Range("F16").Value = 1.050760700887289
$ws.Range("I16").Value = 1.038057503135105
$ws.Range("J16").Value = 1.039178245875094
$ws.Range("K16").Value = 1.039508981561743
$ws.Range("L16").Value = 1.044954781339993
$ws.Range("M16").Value = 1.054146420579465
$ws.Range("N16").Value = 1.040653998179647

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033183305897251
$ws.Range("D17").Value = 1.036324433705192
$ws.Range("E17").Value = 1.04186978541565
$ws.Range("F17").Value = 1.051189975350032
$ws.Range("I17").Value = 1.038146988262087
$ws.Range("J17").Value = 1.039383884099785
$ws.Range("K17").Value = 1.039689143754845
$ws.Range("L17").Value = 1.045215296242382
$ws.Range("M17").Value = 1.054503719200021
$ws.Range("N17").Value = 1.040859928434204

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033383716924842
$ws.Range("D18").Value = 1.036471855787329
$ws.Range("E18").Value = 1.042064047439999
$ws.Range("F18").Value = 1.051440512574222
$ws.Range("I18").Value = 1.038199037320507
$ws.Range("J18").Value = 1.039503791240536
$ws.Range("K18").Value = 1.039794168273754
$ws.Range("L18").Value = 1.045367261280854
$ws.Range("M18").Value = 1.054712204165157
$ws.Range("N18").Value = 1.040980005856849

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033452058375532
$ws.Range("D19").Value = 1.036522126448926
$ws.Range("E19").Value = 1.042130300600704
$ws.Range("F19").Value = 1.051525964484092
$ws.Range("I19").Value = 1.03821675990509
$ws.Range("J19").Value = 1.039544669977049
$ws.Range("K19").Value = 1.039829968538395
$ws.Range("L19").Value = 1.045419079285925
$ws.Range("M19").Value = 1.054783305474281
$ws.Range("N19").Value = 1.041020942645857

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033146444874282
$ws.Range("D20").Value = 1.036297318212905
$ws.Range("E20").Value = 1.041834059425358
$ws.Range("F20").Value = 1.051143902875617
$ws.Range("I20").Value = 1.038137402462634
$ws.Range("J20").Value = 1.039361825004514
$ws.Range("K20").Value = 1.0396698203634
$ws.Range("L20").Value = 1.045187344303385
$ws.Range("M20").Value = 1.054465376289589
$ws.Range("N20").Value = 1.040837838012488

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032153318776869
$ws.Range("D21").Value = 1.035566696415295
$ws.Range("E21").Value = 1.040871983740075
$ws.Range("F21").Value = 1.049903533496875
$ws.Range("I21").Value = 1.037877687888831
$ws.Range("J21").Value = 1.038766931914865
$ws.Range("K21").Value = 1.039148447870482
$ws.Range("L21").Value = 1.044434084140224
$ws.Range("M21").Value = 1.05343268521978
$ws.Range("N21").Value = 1.040242100106416

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031529253690209
$ws.Range("D22").Value = 1.035107525516521
$ws.Range("E22").Value = 1.040267886738376
$ws.Range("F22").Value = 1.049125010418842
$ws.Range("I22").Value = 1.037713082832863
$ws.Range("J22").Value = 1.038392561556793
$ws.Range("K22").Value = 1.038820096885754
$ws.Range("L22").Value = 1.043960586238528
$ws.Range("M22").Value = 1.052784107243864
$ws.Range("N22").Value = 1.039867198099496

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.031860048609637
$ws.Range("D23").Value = 1.03535092126978
$ws.Range("E23").Value = 1.040588053920672
$ws.Range("F23").Value = 1.049537592279922
$ws.Range("I23").Value = 1.037800467272504
$ws.Range("J23").Value = 1.038591054049185
$ws.Range("K23").Value = 1.038994213148505
$ws.Range("L23").Value = 1.04421158600315
$ws.Range("M23").Value = 1.053127862518842
$ws.Range("N23").Value = 1.040065972473994

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033163100665539
$ws.Range("D24").Value = 1.036309570471563
$ws.Range("E24").Value = 1.041850202191837
$ws.Range("F24").Value = 1.051164720584258
$ws.Range("I24").Value = 1.038141734324272
$ws.Range("J24").Value = 1.039371792678124
$ws.Range("K24").Value = 1.039678551960924
$ws.Range("L24").Value = 1.045199974547214
$ws.Range("M24").Value = 1.054482701557177
$ws.Range("N24").Value = 1.040847819841338

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034676719861966
$ws.Range("D25").Value = 1.037422853291998
$ws.Range("E25").Value = 1.043318289107698
$ws.Range("F25").Value = 1.053058741743383
$ws.Range("I25").Value = 1.038532022308103
$ws.Range("J25").Value = 1.04027629626157
$ws.Range("K25").Value = 1.040470291127944
$ws.Range("L25").Value = 1.046347376499715
$ws.Range("M25").Value = 1.056058001550839
$ws.Range("N25").Value = 1.04175360792362

